$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto snapshot values.
# D-column values are forced to remain plain text (matching their original inline-string
# representation) by temporarily applying a text number format, then restoring the
# cell's original style so no residual formatting change is left behind.

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.826.64"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  -2.56%  "

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.094.50"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  -0.80%  "

$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  -0.30%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.54"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.70%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -0.20%  "

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5177"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -1.33%  "

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4467"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -0.76%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09468"
$ws.Range("D9").Style = $origStyle

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.93"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -2.52%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.171"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.34%  "

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.04"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +2.86%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.104.16"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -0.27%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.726"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -0.68%  "

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.035"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -0.25%  "

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.07"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +1.31%  "

$ws.Range("E17").Value = "  +0.22%  "

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.011"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -0.34%  "

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06706"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +0.11%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.55"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +6.50%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.010"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -0.20%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.169"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -2.33%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.889.34"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -2.67%  "

$ws.Range("E24").Value = "  -1.00%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.322"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -2.69%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.348.76"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -0.48%  "

$ws.Range("E27").Value = "  -1.48%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.77"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -0.85%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.525"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.25%  "

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.44"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -1.38%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.156"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -2.99%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1057"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -1.40%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.614"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -1.39%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.214"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -2.19%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.957"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +0.27%  "

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.141"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +4.54%  "

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.13"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -1.17%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02562"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -3.25%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06740"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -1.09%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2274"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -1.64%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.45"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -0.98%  "

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6879"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +0.27%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.297"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +3.14%  "

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6639"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +3.38%  "

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.06"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -5.54%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.278"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -1.23%  "

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.644"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -1.54%  "

$ws.Range("E48").Value = "  -2.75%  "

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000337"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -8.63%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.60"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -1.28%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07133"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -2.18%  "
